# Auto-generated edit script: updates cryptos price/volume table
# to reflect the GitHub-Actions scraper refresh described in the commit
# message ("Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.302.20'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '1.659.36'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('D4').Value = '''1.006'
$ws.Range('E4').Value = '  +0.56%  '
$ws.Range('D5').Value = '''219.52'
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('D6').Value = '''0.5230'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').Value = '''1.006'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').Value = '''0.2665'
$ws.Range('E8').Value = '  +1.71%  '
$ws.Range('D9').Value = '''0.06334'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('D10').Value = '''21.47'
$ws.Range('E10').Value = '  +3.59%  '
$ws.Range('D11').Value = '''0.07757'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.682.37'
$ws.Range('E12').Value = '  +2.05%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.447'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').Value = '''0.5494'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '0.0₅8247'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').Value = '''65.16'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').Value = '26.331.43'
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('D18').Value = '''1.006'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '''4.705'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').Value = '''192.18'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = '''10.23'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '''6.232'
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('D23').Value = '''1.009'
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = '''138.87'
$ws.Range('E24').Value = '  -3.25%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = '''0.1257'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').Value = '''7.338'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').Value = '''16.12'
$ws.Range('E27').Value = '  +0.73%  '
$ws.Range('D28').Value = '''1.425'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').Value = '''0.06078'
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').Value = '''1.289'
$ws.Range('E30').Value = '  +2.30%  '
$ws.Range('D31').Value = '''3.557'
$ws.Range('E31').Value = '  +3.74%  '
$ws.Range('D32').Value = '''3.383'
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('D33').Value = '''1.669'
$ws.Range('E33').Value = '  +1.14%  '
$ws.Range('D34').Value = '''0.9923'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = '''2.427'
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').Value = '''2.768'
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '''0.5987'
$ws.Range('E37').Value = '  +6.47%  '
$ws.Range('D38').Value = '''0.01605'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '''5.987'
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('D40').Value = '1.081.66'
$ws.Range('E40').Value = '  +5.42%  '
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').Value = '''1.005'
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('D43').Value = '''100.11'
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('D44').Value = '1.806.11'
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('D45').Value = '0.0₈110'
$ws.Range('E45').Value = '  +2.70%  '
$ws.Range('D46').Value = '''57.65'
$ws.Range('E46').Value = '  +3.41%  '
$ws.Range('D47').Value = '''8.102'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').Value = '''0.9999'
$ws.Range('E48').Value = '  -0.18%  '
$ws.Range('D49').Value = '''0.05201'
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('D50').Value = '''1.472'
$ws.Range('E50').Value = '  +6.30%  '
$ws.Range('D51').Value = '''0.4239'
$ws.Range('E51').Value = '  +0.73%  '
